$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Trening" column header (F1) marking the training-segment split
$ws.Cells.Item(1, 6).Value = "Trening"

# Rewrite the data rows (A2:F13): the raw GPS export now covers two training
# segments (Duza Gra / Mala Gra) instead of the previous single block, so the
# row count grows from 6 to 12 and every row gets a Trening label.
$ws.Cells.Item(2, 2).Value = 504.3
$ws.Cells.Item(2, 3).Value = 10.13
$ws.Cells.Item(2, 4).Value = 1.69507726601192
$ws.Cells.Item(2, 5).Value = "10-15"
$ws.Cells.Item(2, 6).Value = "Duża Gra"

$ws.Cells.Item(3, 2).Value = 622.3
$ws.Cells.Item(3, 3).Value = 11.6
$ws.Cells.Item(3, 4).Value = 1.628245200429644
$ws.Cells.Item(3, 5).Value = "10-15"
$ws.Cells.Item(3, 6).Value = "Duża Gra"

$ws.Cells.Item(4, 2).Value = 679.2
$ws.Cells.Item(4, 3).Value = 10.76
$ws.Cells.Item(4, 4).Value = 1.650316561971393
$ws.Cells.Item(4, 5).Value = "10-15"
$ws.Cells.Item(4, 6).Value = "Duża Gra"

$ws.Cells.Item(5, 2).Value = 504.2
$ws.Cells.Item(5, 3).Value = 9.56
$ws.Cells.Item(5, 4).Value = 1.709262388093131
$ws.Cells.Item(5, 5).Value = "5-10"
$ws.Cells.Item(5, 6).Value = "Duża Gra"

$ws.Cells.Item(6, 2).Value = 559.2
$ws.Cells.Item(6, 3).Value = 9.39
$ws.Cells.Item(6, 4).Value = 1.576774409839084
$ws.Cells.Item(6, 5).Value = "5-10"
$ws.Cells.Item(6, 6).Value = "Duża Gra"

$ws.Cells.Item(7, 2).Value = 679
$ws.Cells.Item(7, 3).Value = 9.43
$ws.Cells.Item(7, 4).Value = 1.591449413980756
$ws.Cells.Item(7, 5).Value = "5-10"
$ws.Cells.Item(7, 6).Value = "Duża Gra"

$ws.Cells.Item(8, 2).Value = 1131.6
$ws.Cells.Item(8, 3).Value = 11.48
$ws.Cells.Item(8, 4).Value = 3.327181679861887
$ws.Cells.Item(8, 5).Value = "10-15"
$ws.Cells.Item(8, 6).Value = "Mała Gra"

$ws.Cells.Item(9, 2).Value = 1344.8
$ws.Cells.Item(9, 3).Value = 11.94
$ws.Cells.Item(9, 4).Value = 3.132092203412737
$ws.Cells.Item(9, 5).Value = "10-15"
$ws.Cells.Item(9, 6).Value = "Mała Gra"

$ws.Cells.Item(10, 2).Value = 1366.5
$ws.Cells.Item(10, 3).Value = 12.52
$ws.Cells.Item(10, 4).Value = 3.309343031474521
$ws.Cells.Item(10, 5).Value = "10-15"
$ws.Cells.Item(10, 6).Value = "Mała Gra"

$ws.Cells.Item(11, 2).Value = 1131.4
$ws.Cells.Item(11, 3).Value = 8.84
$ws.Cells.Item(11, 4).Value = 3.157211099352156
$ws.Cells.Item(11, 5).Value = "5-10"
$ws.Cells.Item(11, 6).Value = "Mała Gra"

$ws.Cells.Item(12, 2).Value = 1344.6
$ws.Cells.Item(12, 3).Value = 9.07
$ws.Cells.Item(12, 4).Value = 2.913132531302316
$ws.Cells.Item(12, 5).Value = "5-10"
$ws.Cells.Item(12, 6).Value = "Mała Gra"

$ws.Cells.Item(13, 2).Value = 1366.3
$ws.Cells.Item(13, 3).Value = 9.98
$ws.Cells.Item(13, 4).Value = 2.879826562745229
$ws.Cells.Item(13, 5).Value = "5-10"
$ws.Cells.Item(13, 6).Value = "Mała Gra"

# Column A holds Excel date-time serials; give the first cell an explicit
# custom date/time format (registering both the lowercase and uppercase
# format codes), then reuse the resulting format for the remaining rows.
$ws.Cells.Item(2, 1).Value = 45684.59148553241
$ws.Cells.Item(2, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(2, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(3, 1).Value = 45684.59285127315
$ws.Cells.Item(3, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(4, 1).Value = 45684.59350983796
$ws.Cells.Item(4, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(5, 1).Value = 45684.591484375
$ws.Cells.Item(5, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(6, 1).Value = 45684.59212094908
$ws.Cells.Item(6, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(7, 1).Value = 45684.59350752315
$ws.Cells.Item(7, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(8, 1).Value = 45684.59874594907
$ws.Cells.Item(8, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(9, 1).Value = 45684.60121354167
$ws.Cells.Item(9, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(10, 1).Value = 45684.60146469907
$ws.Cells.Item(10, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(11, 1).Value = 45684.59874363426
$ws.Cells.Item(11, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(12, 1).Value = 45684.60121122685
$ws.Cells.Item(12, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(13, 1).Value = 45684.60146238426
$ws.Cells.Item(13, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

Write-Output "done"
